$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mean_return, std_return, prob_outperform, n_obs for rows 2-6 (Monday..Friday)

$ws.Range("C2").Value = 0.0007417590480215013
$ws.Range("D2").Value = 0.02791327872616079
$ws.Range("E2").Value = 0.4998952265322383
$ws.Range("F2").Value = 33057

$ws.Range("C3").Value = 0.001425599175844714
$ws.Range("D3").Value = 0.02739924426249444
$ws.Range("E3").Value = 0.4864314425136891
$ws.Range("F3").Value = 33057

$ws.Range("C4").Value = 0.001100451205830323
$ws.Range("D4").Value = 0.02887899612044406
$ws.Range("E4").Value = 0.4826199068227452
$ws.Range("F4").Value = 33058

$ws.Range("C5").Value = 0.0004865181101752538
$ws.Range("D5").Value = 0.02993613248699877
$ws.Range("E5").Value = 0.480782647206938
$ws.Range("F5").Value = 33058

$ws.Range("C6").Value = 0.0008749186405867679
$ws.Range("D6").Value = 0.02839125587014159
$ws.Range("E6").Value = 0.4452828866115085
$ws.Range("F6").Value = 33059
